$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 362.55554
$ws.Range("I9").Value = 407.33334
$ws.Range("K9").Value = 407.33334
$ws.Range("M9").Value = -238.33334
$ws.Range("H12").Value = 313.66666
$ws.Range("I12").Value = 220.5
$ws.Range("K12").Value = 220.5
$ws.Range("M12").Value = -50.5
$ws.Range("H17").Value = 1103.5333
$ws.Range("J17").Value = 1176.28
$ws.Range("L17").Value = 3528.84
$ws.Range("N17").Value = -3864.84
$ws.Range("H18").Value = 940
$ws.Range("I18").Value = 940
$ws.Range("K18").Value = 940
$ws.Range("M18").Value = -656
$ws.Range("H33").Value = 2870077.8
$ws.Range("I33").Value = 5500309
$ws.Range("K33").Value = 5500309
$ws.Range("M33").Value = -5500080
$ws.Range("H113").Value = 2503.25
$ws.Range("J113").Value = 2497
$ws.Range("L113").Value = 2497
$ws.Range("N113").Value = -9005
$ws.Range("H116").Value = 999998.5
$ws.Range("J116").Value = 1000000
$ws.Range("L116").Value = 1000000
$ws.Range("N116").Value = -1006884
$ws.Range("H132").Value = 114795.78
$ws.Range("I132").Value = 146301.72
$ws.Range("J132").Value = 4525
$ws.Range("K132").Value = 438905.16
$ws.Range("L132").Value = 13575
$ws.Range("M132").Value = -436375.16
$ws.Range("N132").Value = -18635
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1751.8823
$ws.Range("I2").Value = 846.5
$ws.Range("K2").Value = 846.5
$ws.Range("M2").Value = -733.5
$ws.Range("H32").Value = 813.08887
$ws.Range("I32").Value = 824.0909
$ws.Range("K32").Value = 824.0909
$ws.Range("M32").Value = -537.0909
$ws.Range("H45").Value = 1333.2222
$ws.Range("I45").Value = 1323.4117
$ws.Range("K45").Value = 1323.4117
$ws.Range("M45").Value = -946.4117000000001
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H61").Value = 4112.1113
$ws.Range("I61").Value = 3996.2
$ws.Range("K61").Value = 3996.2
$ws.Range("M61").Value = -3784.2
$ws.Range("H116").Value = 1751.8823
$ws.Range("I116").Value = 846.5
$ws.Range("K116").Value = 846.5
$ws.Range("M116").Value = 1447.5
$ws.Range("H136").Value = 4112.1113
$ws.Range("I136").Value = 3996.2
$ws.Range("K136").Value = 11988.6
$ws.Range("M136").Value = -9438.599999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1751.8823
$ws.Range("I3").Value = 846.5
$ws.Range("K3").Value = 846.5
$ws.Range("M3").Value = -732.5
$ws.Range("H14").Value = 699.5
$ws.Range("I14").Value = 999
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 999
$ws.Range("L14").Value = 400
$ws.Range("M14").Value = -827
$ws.Range("N14").Value = -744
$ws.Range("H22").Value = 365.86957
$ws.Range("I22").Value = 257.91666
$ws.Range("J22").Value = 483.63635
$ws.Range("K22").Value = 257.91666
$ws.Range("L22").Value = 483.63635
$ws.Range("M22").Value = -84.91665999999998
$ws.Range("N22").Value = -829.63635
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H24").Value = 2327.75
$ws.Range("I24").Value = 2327.75
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 2327.75
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -2092.75
$ws.Range("N24").ClearContents()
$ws.Range("H25").Value = 4047.5
$ws.Range("I25").Value = 3277
$ws.Range("J25").Value = 7900
$ws.Range("K25").Value = 3277
$ws.Range("L25").Value = 7900
$ws.Range("M25").Value = -3042
$ws.Range("N25").Value = -8370
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 5982
$ws.Range("I31").Value = 5982
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5982
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5730
$ws.Range("N31").ClearContents()
$ws.Range("H36").Value = 2592.5
$ws.Range("I36").Value = 2592.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2592.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2058.5
$ws.Range("N36").ClearContents()
$ws.Range("H37").Value = 3879.3333
$ws.Range("I37").Value = 3879.3333
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3879.3333
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3742.3333
$ws.Range("N37").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H107").Value = 1230.875
$ws.Range("I107").Value = 969.8
$ws.Range("K107").Value = 969.8
$ws.Range("M107").Value = 950.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 6974.75
$ws.Range("J13").Value = 7631.6665
$ws.Range("L13").Value = 7631.6665
$ws.Range("N13").Value = -7909.6665
$ws.Range("H22").Value = 1138.35
$ws.Range("I22").Value = 126.46154
$ws.Range("J22").Value = 3017.5715
$ws.Range("K22").Value = 126.46154
$ws.Range("L22").Value = 3017.5715
$ws.Range("M22").Value = 223.53846
$ws.Range("N22").Value = -3717.5715
$ws.Range("H62").Value = 33948.832
$ws.Range("I62").Value = 29898
$ws.Range("K62").Value = 29898
$ws.Range("M62").Value = -29274
$ws.Range("H65").Value = 33948.832
$ws.Range("I65").Value = 29898
$ws.Range("K65").Value = 149490
$ws.Range("M65").Value = -146370
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 140
$ws.Range("I14").Value = 140
$ws.Range("K14").Value = 420
$ws.Range("M14").Value = -247
$ws.Range("H70").Value = 639
$ws.Range("I70").Value = 639
$ws.Range("K70").Value = 1917
$ws.Range("M70").Value = -1602
$ws.Range("H73").Value = 639
$ws.Range("I73").Value = 639
$ws.Range("K73").Value = 1917
$ws.Range("M73").Value = -825
$ws.Range("H103").Value = 120.5
$ws.Range("I103").Value = 129.125
$ws.Range("J103").Value = 103.25
$ws.Range("K103").Value = 387.375
$ws.Range("L103").Value = 309.75
$ws.Range("M103").Value = 491.625
$ws.Range("N103").Value = -2067.75
$ws.Range("H114").Value = 4166.6
$ws.Range("J114").Value = 4166.6
$ws.Range("L114").Value = 12499.8
$ws.Range("N114").Value = -19007.8
$ws.Range("H117").Value = 3393.75
$ws.Range("I117").Value = 2500
$ws.Range("J117").Value = 3691.6667
$ws.Range("K117").Value = 7500
$ws.Range("L117").Value = 11075.0001
$ws.Range("M117").Value = -4058
$ws.Range("N117").Value = -17959.0001
$ws.Range("H131").Value = 1213191.5
$ws.Range("I131").Value = 839.44446
$ws.Range("K131").Value = 2518.33338
$ws.Range("M131").Value = 2521.66662
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9484.77
$ws.Range("I70").Value = 9425.5
$ws.Range("J70").Value = 9579.6
$ws.Range("K70").Value = 9425.5
$ws.Range("L70").Value = 9579.6
$ws.Range("M70").Value = -9155.5
$ws.Range("N70").Value = -10119.6
$ws.Range("H73").Value = 9484.77
$ws.Range("I73").Value = 9425.5
$ws.Range("J73").Value = 9579.6
$ws.Range("K73").Value = 9425.5
$ws.Range("L73").Value = 9579.6
$ws.Range("M73").Value = -8489.5
$ws.Range("N73").Value = -11451.6
$ws.Range("H102").Value = 1498.4445
$ws.Range("I102").Value = 1373.375
$ws.Range("J102").Value = 2499
$ws.Range("K102").Value = 1373.375
$ws.Range("L102").Value = 2499
$ws.Range("M102").Value = 248.625
$ws.Range("N102").Value = -5743
$ws.Range("H132").Value = 2412.0667
$ws.Range("I132").Value = 2253.5557
$ws.Range("K132").Value = 6760.6671
$ws.Range("M132").Value = -4230.6671
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3432.111
$ws.Range("I61").Value = 2127.1428
$ws.Range("J61").Value = 7999.5
$ws.Range("K61").Value = 2127.1428
$ws.Range("L61").Value = 7999.5
$ws.Range("M61").Value = -1925.1428
$ws.Range("N61").Value = -8403.5
$ws.Range("H113").Value = 3432.111
$ws.Range("I113").Value = 2127.1428
$ws.Range("J113").Value = 7999.5
$ws.Range("K113").Value = 2127.1428
$ws.Range("L113").Value = 7999.5
$ws.Range("M113").Value = 42.85719999999992
$ws.Range("N113").Value = -12339.5
$ws.Range("H122").Value = 3397.2666
$ws.Range("I122").Value = 3021
$ws.Range("K122").Value = 9063
$ws.Range("M122").Value = -6613
$ws.Range("H132").Value = 4293.7
$ws.Range("I132").Value = 4421.7856
$ws.Range("K132").Value = 13265.3568
$ws.Range("M132").Value = -10735.3568
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 11798.875
$ws.Range("I52").Value = 5442
$ws.Range("J52").Value = 12707
$ws.Range("K52").Value = 5442
$ws.Range("L52").Value = 12707
$ws.Range("M52").Value = -5216
$ws.Range("N52").Value = -13159
$ws.Range("H92").Value = 99664.664
$ws.Range("J92").Value = 99664.664
$ws.Range("L92").Value = 99664.664
$ws.Range("N92").Value = -104656.664
$ws.Range("H125").Value = 7776
$ws.Range("J125").Value = 7776
$ws.Range("L125").Value = 7776
$ws.Range("N125").Value = -17616
